$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102 (shifts existing rows 102:212 down to 103:213)
$ws.Rows("102:102").Insert()

# Fill in the constant columns (same values repeated throughout the dataset)
$ws.Range("A102").Value = 4
$ws.Range("B102").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C102").Value = "Los Lagos"
$ws.Range("E102").Value = 10
$ws.Range("F102").Value = 100112009
$ws.Range("G102").Value = "Acelga"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("R102").Value = "Hortaliza"

# Fill in the new weekly data point values
$ws.Range("D102").Value = 44810
$ws.Range("J102").Value = 200
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 3000
$ws.Range("N102").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O102").Value = "Región del Maule"
$ws.Range("P102").Value = 750
$ws.Range("Q102").Value = 4
